$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of combined data (Dan / marinescu) matching the
# existing layout used for Omer/Avisror and Liat/Mulian rows.
$ws.Cells.Item(3, 1).Value = "204264543"
$ws.Cells.Item(3, 2).Value = "Dan"
$ws.Cells.Item(3, 3).Value = "marinescu"
$ws.Cells.Item(3, 4).Value = "1234"
$ws.Cells.Item(3, 5).Value = $true
